# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计"),
#    populated with the quarter's fund-holding detail rows.
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet, shifting the
#    existing quarterly summary rows down by one and renumbering the
#    leading index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" sheet (fund holdings detail)
# ---------------------------------------------------------------------
$q4sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4sheet)
$newSheet.Name = "2022-Q1"

# Match the sheet-level outline/page setup used by every other quarter
# sheet in this workbook.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Reuse the existing look (bold/centered/bordered header row, styled
# index column) from the "2021-Q4" sheet instead of hand-building xfs.
$q4sheet.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)
$q4sheet.Range("A2:H2").Copy()
$newSheet.Range("A2:H8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B:G hold text (fund codes keep leading zeros, the rest are
# percentages/sizes stored as text in this workbook) - force text so
# Excel doesn't "helpfully" reinterpret them as numbers.
$textRange = $newSheet.Range("B2:G8")
$textRange.NumberFormat = "@"

$rows = @(
    @(0, "506001", "万家科创板 2 年定期开放混合型证券投资基金", "12.84", "98.14", "3.46", "0.4443", 7),
    @(1, "007012", "湘财长顺混合A", "4.70", "94.08", "6.08", "0.2858", 9),
    @(2, "008128", "湘财长源股票A", "2.74", "94.29", "6.26", "0.1715", 8),
    @(3, "007013", "湘财长顺混合C", "2.47", "94.08", "6.08", "0.1502", 9),
    @(4, "008129", "湘财长源股票C", "1.05", "94.29", "6.26", "0.0657", 8),
    @(5, "970020", "信达价值精选一年持有期灵活配置混合A", "0.64", "56.02", "4.88", "0.0312", 6),
    @(6, "970021", "信达价值精选一年持有期灵活配置混合B", "0.53", "56.02", "4.88", "0.0259", 6)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Drop the temporary "@" number format again (its only purpose was to
# stop the auto numeric conversion above) without disturbing the text
# that's now stored in those cells.
$blank = $newSheet.Range("Z100")
$blank.Copy()
$textRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Part 2: "总计" sheet - insert the 2022-Q1 summary row at the top
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Push existing data rows (2-6) down to (3-7), working from the bottom
# up so nothing gets clobbered before it's copied.
$totalSheet.Range("A6:D6").Copy($totalSheet.Range("A7:D7"))
$totalSheet.Range("A5:D5").Copy($totalSheet.Range("A6:D6"))
$totalSheet.Range("A4:D4").Copy($totalSheet.Range("A5:D5"))
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

# New first data row.
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 1.17

# Renumber the leading index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# Leave the view the way it started (first sheet focused).
$wb.Worksheets.Item(1).Activate()
